# Fruta / hortaliza, semanal
# Insert a new weekly record as row 64, shifting the existing rows 64-89 down to 65-90.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 64 (pushes old rows 64..89 down to 65..90)
$ws.Rows.Item(64).Insert()

# Populate the newly inserted row 64 with the new record's data
$ws.Cells.Item(64, 1).Value = 2
$ws.Cells.Item(64, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(64, 3).Value = "Coquimbo"
$ws.Cells.Item(64, 4).Value = 44510
$ws.Cells.Item(64, 5).Value = 4
$ws.Cells.Item(64, 6).Value = 100112043
$ws.Cells.Item(64, 7).Value = "Pepino ensalada"
$ws.Cells.Item(64, 8).Value = "Sin especificar"
$ws.Cells.Item(64, 9).Value = "Primera"
$ws.Cells.Item(64, 10).Value = 400
$ws.Cells.Item(64, 11).Value = 5500
$ws.Cells.Item(64, 12).Value = 6000
$ws.Cells.Item(64, 13).Value = 5750
$ws.Cells.Item(64, 14).Value = "$/caja 70 unidades"
$ws.Cells.Item(64, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(64, 16).Value = 82
$ws.Cells.Item(64, 17).Value = 70
$ws.Cells.Item(64, 18).Value = "Hortaliza"
